$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.086.81"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.928.55"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.80"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.18"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.561"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.50"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("E11").Value = "  +3.56%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.02"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.75"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "3.390.43"
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "2.919.26"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.985"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "52.124.98"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.28"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.05"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.77"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +8.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.07"
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.96"
$ws.Range("E29").Value = "  +12.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.61"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.102"
$ws.Range("E31").Value = "  +12.51%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.04"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.03"
$ws.Range("E34").Value = "  +5.66%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.07"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +4.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.61"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.90"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.191.94"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.51"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.29"
$ws.Range("E48").Value = "  -8.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.248"
$ws.Range("E49").Value = "  +10.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0345"
$ws.Range("E50").Value = "  +6.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.951"
$ws.Range("E51").Value = "  -7.64%  "

Write-Host "Update complete"
